$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Staging cell used to force text-typed values (via leading apostrophe / quote-prefix),
# then Paste-Special-Values copies the TEXT into the target cell without carrying any
# cell style along (matching the source workbook, where these numeric-looking strings
# are stored as shared-string text with no explicit "s" style attribute).
$stage = $ws.Cells.Item(5000, 700)

$stage.Value = "'101.1"
$stage.Copy()
$ws.Range("AI2").PasteSpecial(-4163)
$stage.Value = "'112.1"
$stage.Copy()
$ws.Range("AI10").PasteSpecial(-4163)
$stage.Value = "'110.1"
$stage.Copy()
$ws.Range("AI11").PasteSpecial(-4163)
$stage.Value = "'105.49"
$stage.Copy()
$ws.Range("AI12").PasteSpecial(-4163)
$stage.Value = "'103.78"
$stage.Copy()
$ws.Range("AI13").PasteSpecial(-4163)
$stage.Value = "'59125.00"
$stage.Copy()
$ws.Range("AI15").PasteSpecial(-4163)
$stage.Value = "'118.0"
$stage.Copy()
$ws.Range("AI19").PasteSpecial(-4163)
$stage.Value = "'116.3"
$stage.Copy()
$ws.Range("AI27").PasteSpecial(-4163)
$stage.Value = "'112.8"
$stage.Copy()
$ws.Range("AI28").PasteSpecial(-4163)
$stage.Value = "'106.70"
$stage.Copy()
$ws.Range("AI29").PasteSpecial(-4163)
$stage.Value = "'106.61"
$stage.Copy()
$ws.Range("AH30").PasteSpecial(-4163)
$stage.Value = "'108.91"
$stage.Copy()
$ws.Range("AI30").PasteSpecial(-4163)
$stage.Value = "'37873.67"
$stage.Copy()
$ws.Range("AI32").PasteSpecial(-4163)
$stage.Value = "'130.8"
$stage.Copy()
$ws.Range("AI36").PasteSpecial(-4163)
$stage.Value = "'129.8"
$stage.Copy()
$ws.Range("AI44").PasteSpecial(-4163)
$stage.Value = "'261.6"
$stage.Copy()
$ws.Range("AI45").PasteSpecial(-4163)
$stage.Value = "'105.90"
$stage.Copy()
$ws.Range("AI46").PasteSpecial(-4163)
$stage.Value = "'246"
$stage.Copy()
$ws.Range("AI47").PasteSpecial(-4163)
$stage.Value = "'26060.68"
$stage.Copy()
$ws.Range("AI49").PasteSpecial(-4163)
$stage.Value = "'82.2"
$stage.Copy()
$ws.Range("AI53").PasteSpecial(-4163)
$stage.Value = "'110.6"
$stage.Copy()
$ws.Range("AI61").PasteSpecial(-4163)
$stage.Value = "'99.2"
$stage.Copy()
$ws.Range("AI62").PasteSpecial(-4163)
$stage.Value = "'97.40"
$stage.Copy()
$ws.Range("AI63").PasteSpecial(-4163)
$stage.Value = "'92.65"
$stage.Copy()
$ws.Range("AI64").PasteSpecial(-4163)
$stage.Value = "'42948.67"
$stage.Copy()
$ws.Range("AI66").PasteSpecial(-4163)
$stage.Value = "'94.2"
$stage.Copy()
$ws.Range("AI70").PasteSpecial(-4163)
$stage.Value = "'103.8"
$stage.Copy()
$ws.Range("AI78").PasteSpecial(-4163)
$stage.Value = "'105.0"
$stage.Copy()
$ws.Range("AI79").PasteSpecial(-4163)
$stage.Value = "'106.41"
$stage.Copy()
$ws.Range("AI80").PasteSpecial(-4163)
$stage.Value = "'95.03"
$stage.Copy()
$ws.Range("AI81").PasteSpecial(-4163)
$stage.Value = "'21821.34"
$stage.Copy()
$ws.Range("AI83").PasteSpecial(-4163)
$stage.Value = "'99.9"
$stage.Copy()
$ws.Range("AI87").PasteSpecial(-4163)
$stage.Value = "'108.3"
$stage.Copy()
$ws.Range("AI95").PasteSpecial(-4163)
$stage.Value = "'122.0"
$stage.Copy()
$ws.Range("AI96").PasteSpecial(-4163)
$stage.Value = "'103.66"
$stage.Copy()
$ws.Range("AI97").PasteSpecial(-4163)
$stage.Value = "'105.45"
$stage.Copy()
$ws.Range("AI98").PasteSpecial(-4163)
$stage.Value = "'37719.04"
$stage.Copy()
$ws.Range("AI100").PasteSpecial(-4163)
$stage.Value = "'101.0"
$stage.Copy()
$ws.Range("AI104").PasteSpecial(-4163)
$stage.Value = "'101.9"
$stage.Copy()
$ws.Range("AI112").PasteSpecial(-4163)
$stage.Value = "'112.9"
$stage.Copy()
$ws.Range("AI113").PasteSpecial(-4163)
$stage.Value = "'100.00"
$stage.Copy()
$ws.Range("AI114").PasteSpecial(-4163)
$stage.Value = "'98.06"
$stage.Copy()
$ws.Range("AI115").PasteSpecial(-4163)
$stage.Value = "'24103.08"
$stage.Copy()
$ws.Range("AI117").PasteSpecial(-4163)
$stage.Value = "'109.8"
$stage.Copy()
$ws.Range("AI121").PasteSpecial(-4163)
$stage.Value = "'122.0"
$stage.Copy()
$ws.Range("AI129").PasteSpecial(-4163)
$stage.Value = "'118.9"
$stage.Copy()
$ws.Range("AI130").PasteSpecial(-4163)
$stage.Value = "'105.11"
$stage.Copy()
$ws.Range("AI131").PasteSpecial(-4163)
$stage.Value = "'108.18"
$stage.Copy()
$ws.Range("AI132").PasteSpecial(-4163)
$stage.Value = "'304091.34"
$stage.Copy()
$ws.Range("AI134").PasteSpecial(-4163)
$stage.Value = "'99.6"
$stage.Copy()
$ws.Range("AI138").PasteSpecial(-4163)
$stage.Value = "'104.2"
$stage.Copy()
$ws.Range("AI146").PasteSpecial(-4163)
$stage.Value = "'94.9"
$stage.Copy()
$ws.Range("AI147").PasteSpecial(-4163)
$stage.Value = "'96.90"
$stage.Copy()
$ws.Range("AI148").PasteSpecial(-4163)
$stage.Value = "'91.08"
$stage.Copy()
$ws.Range("AI149").PasteSpecial(-4163)
$stage.Value = "'19935.44"
$stage.Copy()
$ws.Range("AI151").PasteSpecial(-4163)
$stage.Value = "'69.6"
$stage.Copy()
$ws.Range("AI155").PasteSpecial(-4163)
$stage.Value = "'105.6"
$stage.Copy()
$ws.Range("AI163").PasteSpecial(-4163)
$stage.Value = "'102.5"
$stage.Copy()
$ws.Range("AI164").PasteSpecial(-4163)
$stage.Value = "'122.72"
$stage.Copy()
$ws.Range("AI165").PasteSpecial(-4163)
$stage.Value = "'108.34"
$stage.Copy()
$ws.Range("AI166").PasteSpecial(-4163)
$stage.Value = "'128934.59"
$stage.Copy()
$ws.Range("AI168").PasteSpecial(-4163)
$stage.Value = "'105.3"
$stage.Copy()
$ws.Range("AI172").PasteSpecial(-4163)
$stage.Value = "'111.7"
$stage.Copy()
$ws.Range("AI180").PasteSpecial(-4163)
$stage.Value = "'138.9"
$stage.Copy()
$ws.Range("AI181").PasteSpecial(-4163)
$stage.Value = "'103.62"
$stage.Copy()
$ws.Range("AI182").PasteSpecial(-4163)
$stage.Value = "'110.85"
$stage.Copy()
$ws.Range("AI183").PasteSpecial(-4163)
$stage.Value = "'15112.07"
$stage.Copy()
$ws.Range("AI185").PasteSpecial(-4163)
$stage.Value = "'104.7"
$stage.Copy()
$ws.Range("AI189").PasteSpecial(-4163)
$stage.Value = "'108.7"
$stage.Copy()
$ws.Range("AI197").PasteSpecial(-4163)
$stage.Value = "'107.2"
$stage.Copy()
$ws.Range("AI198").PasteSpecial(-4163)
$stage.Value = "'112.68"
$stage.Copy()
$ws.Range("AI199").PasteSpecial(-4163)
$stage.Value = "'108.09"
$stage.Copy()
$ws.Range("AI200").PasteSpecial(-4163)
$stage.Value = "'46325.47"
$stage.Copy()
$ws.Range("AI202").PasteSpecial(-4163)
$stage.Value = "'104.9"
$stage.Copy()
$ws.Range("AI206").PasteSpecial(-4163)
$stage.Value = "'110.2"
$stage.Copy()
$ws.Range("AI214").PasteSpecial(-4163)
$stage.Value = "'104.4"
$stage.Copy()
$ws.Range("AI215").PasteSpecial(-4163)
$stage.Value = "'107.23"
$stage.Copy()
$ws.Range("AI216").PasteSpecial(-4163)
$stage.Value = "'107.84"
$stage.Copy()
$ws.Range("AI217").PasteSpecial(-4163)
$stage.Value = "'38682.54"
$stage.Copy()
$ws.Range("AI219").PasteSpecial(-4163)

$ws.Application.CutCopyMode = $false
$stage.Clear()

# New AutoFilter over the full data range (header row 1 through the last data row 222,
# columns A through AI) -- matches the sheet's existing hidden _FilterDatabase name.
$ws.Range("A1:AI222").AutoFilter()

# Restore the active selection to X19 (as in the committed workbook).
$ws.Range("X19").Select()

